$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (class "0.0")
$ws.Range("B2").Value = 0.9815724815724816
$ws.Range("C2").Value = 0.8619201725997843
$ws.Range("D2").Value = 0.9178632969557725
$ws.Range("E2").Value = 927

# Row 3 (class "1.0")
$ws.Range("B3").Value = 0.2849162011173184
$ws.Range("D3").Value = 0.4163265306122448

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.8559919436052367
$ws.Range("C4").Value = 0.8559919436052367
$ws.Range("D4").Value = 0.8559919436052367
$ws.Range("E4").Value = 0.8559919436052367

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.6332443413449
$ws.Range("C5").Value = 0.8173237226635285
$ws.Range("D5").Value = 0.6670949137840086
$ws.Range("E5").Value = 993

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9352690429923801
$ws.Range("C6").Value = 0.8559919436052367
$ws.Range("D6").Value = 0.8845285269873205
$ws.Range("E6").Value = 993

$wb.Save()
